# Auto-generated script applying scheduled-runner profit recalculation updates
# to the Tonberry_Profits workbook. Updates currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H-N) across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2101.3572
$ws.Range("I19").Value = 979.25
$ws.Range("J19").Value = 2550.2
$ws.Range("K19").Value = 979.25
$ws.Range("L19").Value = 2550.2
$ws.Range("M19").Value = -804.25
$ws.Range("N19").Value = -2900.2
$ws.Range("H32").Value = 3091
$ws.Range("J32").Value = 2909.2
$ws.Range("L32").Value = 2909.2
$ws.Range("N32").Value = -3561.2
$ws.Range("H43").Value = 1491.8334
$ws.Range("I43").Value = 1549.4
$ws.Range("J43").Value = 1450.7142
$ws.Range("K43").Value = 1549.4
$ws.Range("L43").Value = 1450.7142
$ws.Range("M43").Value = -1480.4
$ws.Range("N43").Value = -1588.7142
$ws.Range("H94").Value = 2821.8
$ws.Range("I94").Value = 2821.8
$ws.Range("K94").Value = 2821.8
$ws.Range("M94").Value = -2370.8
$ws.Range("H106").Value = 3144.182
$ws.Range("I106").Value = 4243.857
$ws.Range("J106").Value = 1219.75
$ws.Range("K106").Value = 4243.857
$ws.Range("L106").Value = 1219.75
$ws.Range("M106").Value = -3612.857
$ws.Range("N106").Value = -2481.75
$ws.Range("H113").Value = 20029.705
$ws.Range("I113").Value = 34756.445
$ws.Range("J113").Value = 3462.125
$ws.Range("K113").Value = 34756.445
$ws.Range("L113").Value = 3462.125
$ws.Range("M113").Value = -31502.445
$ws.Range("N113").Value = -9970.125
$ws.Range("H138").Value = 2530.2769
$ws.Range("I138").Value = 3292.5417
$ws.Range("J138").Value = 2084.0732
$ws.Range("K138").Value = 9877.625100000001
$ws.Range("L138").Value = 6252.219599999999
$ws.Range("M138").Value = -4737.625100000001
$ws.Range("N138").Value = -16532.2196
$ws.Range("H139").Value = 67154.5
$ws.Range("J139").Value = 73600
$ws.Range("L139").Value = 73600
$ws.Range("N139").Value = -83880
$ws.Range("H140").Value = 83229.3
$ws.Range("J140").Value = 83229.3
$ws.Range("L140").Value = 83229.3
$ws.Range("N140").Value = -93589.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5042.15
$ws.Range("I32").Value = 3843.3333
$ws.Range("K32").Value = 3843.3333
$ws.Range("M32").Value = -3556.3333
$ws.Range("H61").Value = 5357.273
$ws.Range("I61").Value = 4913.32
$ws.Range("J61").Value = 6744.625
$ws.Range("K61").Value = 4913.32
$ws.Range("L61").Value = 6744.625
$ws.Range("M61").Value = -4701.32
$ws.Range("N61").Value = -7168.625
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H130").Value = 50676.363
$ws.Range("J130").Value = 50676.363
$ws.Range("L130").Value = 50676.363
$ws.Range("N130").Value = -60716.363
$ws.Range("H132").Value = 1581.5405
$ws.Range("I132").Value = 1194.8276
$ws.Range("J132").Value = 2983.375
$ws.Range("K132").Value = 3584.4828
$ws.Range("L132").Value = 8950.125
$ws.Range("M132").Value = -1054.4828
$ws.Range("N132").Value = -14010.125
$ws.Range("H136").Value = 5357.273
$ws.Range("I136").Value = 4913.32
$ws.Range("J136").Value = 6744.625
$ws.Range("K136").Value = 14739.96
$ws.Range("L136").Value = 20233.875
$ws.Range("M136").Value = -12189.96
$ws.Range("N136").Value = -25333.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2442.6155
$ws.Range("I20").Value = 2214.6667
$ws.Range("J20").Value = 2638
$ws.Range("K20").Value = 2214.6667
$ws.Range("L20").Value = 2638
$ws.Range("M20").Value = -1967.6667
$ws.Range("N20").Value = -3132
$ws.Range("H80").Value = 5154.4546
$ws.Range("J80").Value = 5659.9
$ws.Range("L80").Value = 5659.9
$ws.Range("N80").Value = -7655.9
$ws.Range("H83").Value = 5154.4546
$ws.Range("J83").Value = 5659.9
$ws.Range("L83").Value = 28299.5
$ws.Range("N83").Value = -38283.5
$ws.Range("H134").Value = 4193.5127
$ws.Range("I134").Value = 4348.243
$ws.Range("J134").Value = 1331
$ws.Range("K134").Value = 13044.729
$ws.Range("L134").Value = 3993
$ws.Range("M134").Value = -10509.729
$ws.Range("N134").Value = -9063

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 680.13336
$ws.Range("I16").Value = 599.36365
$ws.Range("K16").Value = 599.36365
$ws.Range("M16").Value = -312.36365
$ws.Range("H31").Value = 3138.375
$ws.Range("I31").Value = 3477.2
$ws.Range("K31").Value = 3477.2
$ws.Range("M31").Value = -3182.2
$ws.Range("H34").Value = 3138.375
$ws.Range("I34").Value = 3477.2
$ws.Range("K34").Value = 3477.2
$ws.Range("M34").Value = -3275.2
$ws.Range("H62").Value = 2906
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 2906
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H113").Value = 680.13336
$ws.Range("I113").Value = 599.36365
$ws.Range("K113").Value = 599.36365
$ws.Range("M113").Value = 1570.63635
$ws.Range("H132").Value = 1812.6207
$ws.Range("I132").Value = 1005.4545
$ws.Range("J132").Value = 4349.4287
$ws.Range("K132").Value = 3016.3635
$ws.Range("L132").Value = 13048.2861
$ws.Range("M132").Value = -486.3635000000004
$ws.Range("N132").Value = -18108.2861
$ws.Range("H134").Value = 1826.9524
$ws.Range("I134").Value = 1680.2222
$ws.Range("K134").Value = 5040.6666
$ws.Range("M134").Value = -2505.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 216.5
$ws.Range("I6").Value = 299.75
$ws.Range("K6").Value = 899.25
$ws.Range("M6").Value = -786.25
$ws.Range("H11").Value = 827.25
$ws.Range("I11").Value = 770
$ws.Range("K11").Value = 2310
$ws.Range("M11").Value = -2170
$ws.Range("H33").Value = 121.3
$ws.Range("I33").Value = 134
$ws.Range("J33").Value = 102.25
$ws.Range("K33").Value = 804
$ws.Range("L33").Value = 613.5
$ws.Range("M33").Value = -521
$ws.Range("N33").Value = -1179.5
$ws.Range("H68").Value = 855.1
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 855.1
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 2565.3
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -4187.3
$ws.Range("H71").Value = 855.1
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 855.1
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 7695.900000000001
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -15807.9
$ws.Range("H122").Value = 767.9286
$ws.Range("I122").Value = 635
$ws.Range("J122").Value = 867.625
$ws.Range("K122").Value = 5715
$ws.Range("L122").Value = 7808.625
$ws.Range("M122").Value = -3265
$ws.Range("N122").Value = -12708.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 15100
$ws.Range("J26").Value = 15100
$ws.Range("L26").Value = 15100
$ws.Range("N26").Value = -15660
$ws.Range("H50").Value = 15100
$ws.Range("J50").Value = 15100
$ws.Range("L50").Value = 15100
$ws.Range("N50").Value = -16096
$ws.Range("H97").Value = 912.4375
$ws.Range("I97").Value = 926.125
$ws.Range("J97").Value = 871.375
$ws.Range("K97").Value = 926.125
$ws.Range("L97").Value = 871.375
$ws.Range("M97").Value = -430.125
$ws.Range("N97").Value = -1863.375
$ws.Range("H113").Value = 805.1739
$ws.Range("I113").Value = 569.4286
$ws.Range("J113").Value = 1171.8889
$ws.Range("K113").Value = 569.4286
$ws.Range("L113").Value = 1171.8889
$ws.Range("M113").Value = 1600.5714
$ws.Range("N113").Value = -5511.8889
$ws.Range("H132").Value = 1242883.9
$ws.Range("I132").Value = 1749948.5
$ws.Range("K132").Value = 5249845.5
$ws.Range("M132").Value = -5247315.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2580.2368
$ws.Range("I132").Value = 1494.0667
$ws.Range("K132").Value = 4482.2001
$ws.Range("M132").Value = -1952.2001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1300.8
$ws.Range("I100").Value = 1126.25
$ws.Range("J100").Value = 1999
$ws.Range("K100").Value = 2252.5
$ws.Range("L100").Value = 3998
$ws.Range("M100").Value = -1711.5
$ws.Range("N100").Value = -5080
$ws.Range("H113").Value = 537.5833
$ws.Range("I113").Value = 371.7647
$ws.Range("J113").Value = 940.2857
$ws.Range("K113").Value = 1115.2941
$ws.Range("L113").Value = 2820.8571
$ws.Range("M113").Value = 1054.7059
$ws.Range("N113").Value = -7160.8571
$ws.Range("H132").Value = 1462.4857
$ws.Range("I132").Value = 1032.8276
$ws.Range("K132").Value = 3098.4828
$ws.Range("M132").Value = -568.4828000000002
$ws.Range("H136").Value = 15874280
$ws.Range("I136").Value = 21368436
$ws.Range("K136").Value = 64105308
$ws.Range("M136").Value = -64102758
